$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: rename the old "Obj2.2" row (currently row 4) to "Obj2.2.1" ---
$ws.Range("A4").Value = "Obj2.2.1"

# --- Step 2: insert a new row right after it for "Obj2.2.2" ---
$ws.Rows.Item(5).Insert()
$ws.Range("A5").Value = "Obj2.2.2"
$ws.Range("B5").Value = "Lambda in PopSim; lambda in model"
$ws.Range("C5").Value = "all three"
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = "variable"
$ws.Range("F5").Value = "iteratively fixed"
$ws.Range("G5").Value = "Uniform 0.5-0.95"
$ws.Range("H5").Value = "base case + lambda"

# --- Step 3: insert a new row after "Obj1.1" (row 2) for "Obj1.2" ---
$ws.Rows.Item(3).Insert()
$ws.Range("A3").Value = "Obj1.2"
$ws.Range("B3").Value = "Initial model validation w/ uninformed priors"
$ws.Range("C3").Value = "all three"
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = "NA"
$ws.Range("G3").Value = "Uniform 0.5-0.95"
$ws.Range("H3").Value = "base-case (equation from original CKMR paper)"
$ws.Range("I3").Value = 'In this case, the "truth" will be the mean abundance over the cohort years.'
$ws.Range("I3").WrapText = $true

# --- Step 4: update the survival prior for "Obj1.1" (row 2) ---
$ws.Range("G2").Value = "beta w/ 10% CV"

# Rows shifted by the inserts above keep their original content:
#  old row 3 (Obj2.1)            -> row 4
#  old row 5 (Obj2.3?)           -> row 7
#  old row 6 (Obj3.1)            -> row 8
#  old row 7 (Obj3.2)            -> row 9
#  old row 8 (Obj4.1)            -> row 10
#  old rows 21-23 (supplemental) -> rows 23-25

# --- Selection, to match the saved workbook view ---
$ws.Range("B16").Select()
